# Update the four pressure/measurement readings in column B (rows 2-5)
# from 6.85 to 5.15, as in the source commit. Switch to manual
# calculation first so dependent formula cells in column D (which the
# diff does not touch) are not force-recalculated by this script -
# only the raw input values change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$excel.Calculation = -4135   # xlCalculationManual

$ws.Range("B2").Value = 5.15
$ws.Range("B3").Value = 5.15
$ws.Range("B4").Value = 5.15
$ws.Range("B5").Value = 5.15
